$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.575.23'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.159.91'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.48'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.15'
$ws.Range('E6').Value = '  -1.56%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.159.99'
$ws.Range('E8').Value = '  +0.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.42'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.471'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000258'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.63'
$ws.Range('E14').Value = '  -1.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.680.00'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.590.26'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.159.48'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.87'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.50'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.61'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('E22').Value = '  +2.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.97'
$ws.Range('E23').Value = '  +3.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.82'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.93'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -3.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.65'
$ws.Range('E28').Value = '  +2.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.16'
$ws.Range('E30').Value = '  -2.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.11'
$ws.Range('E31').Value = '  -5.03%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.67'
$ws.Range('E33').Value = '  -1.15%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.49'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  +1.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0785'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.00'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.20'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.19'
$ws.Range('E39').Value = '  +3.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '462.52'
$ws.Range('E40').Value = '  +2.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0399'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('E42').Value = '  -3.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.33'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.863.43'
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E45').Value = '  +2.78%  '
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.48'
$ws.Range('E47').Value = '  +7.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.57'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.77'
$ws.Range('E50').Value = '  +9.59%  '
$ws.Range('E51').Value = '  -0.33%  '
